$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation is inserted at row 33, pushing the
# existing rows 33:68 down to 34:69 (dates/prices/origins unchanged).
$ws.Rows("33:33").Insert()

$ws.Range("A33").Value = 11
$ws.Range("B33").Value = "Vega Monumental Concepción"
$ws.Range("C33").Value = "Bíobío"
$ws.Range("D33").Value = 44467
$ws.Range("E33").Value = 8
$ws.Range("F33").Value = "Fruta"
$ws.Range("G33").Value = 100108
$ws.Range("H33").Value = "Tropicales y subtropicales"
$ws.Range("I33").Value = 100108002
$ws.Range("J33").Value = "Mango"
$ws.Range("K33").Value = "Sin especificar"
$ws.Range("L33").Value = "Primera"
$ws.Range("M33").Value = 300
$ws.Range("N33").Value = 7500
$ws.Range("O33").Value = 8000
$ws.Range("P33").Value = 7667
$ws.Range("Q33").Value = "$/bandeja 4 kilos"
$ws.Range("R33").Value = "Brasil"
$ws.Range("S33").Value = 1917
$ws.Range("T33").Value = 4
